$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 137.03448
$ws.Range("I33").Value = 97.708336
$ws.Range("J33").Value = 325.8
$ws.Range("K33").Value = 97.708336
$ws.Range("L33").Value = 325.8
$ws.Range("M33").Value = 131.291664
$ws.Range("N33").Value = -783.8

$ws.Range("H106").Value = 3560.4285
$ws.Range("I106").Value = 3560.4285
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3560.4285
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2929.4285

$ws.Range("H116").Value = 3313.2778
$ws.Range("I116").Value = 2639.7273
$ws.Range("J116").Value = 4371.7144
$ws.Range("K116").Value = 2639.7273
$ws.Range("L116").Value = 4371.7144
$ws.Range("M116").Value = 802.2727
$ws.Range("N116").Value = -11255.7144

$ws.Range("H129").Value = 913.5135
$ws.Range("I129").Value = 750
$ws.Range("J129").Value = 918.05554
$ws.Range("K129").Value = 2250
$ws.Range("L129").Value = 2754.16662
$ws.Range("M129").Value = 2750
$ws.Range("N129").Value = -12754.16662

$ws.Range("H132").Value = 5954362
$ws.Range("I132").Value = 7938494
$ws.Range("J132").Value = 1965.7858
$ws.Range("K132").Value = 23815482
$ws.Range("L132").Value = 5897.357400000001
$ws.Range("M132").Value = -23812952
$ws.Range("N132").Value = -10957.3574

$ws.Range("H135").Value = 204.77777
$ws.Range("I135").Value = 204.77777
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 1842.99993
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = 692.0000700000001

$ws.Range("H137").Value = 1453.8
$ws.Range("I137").Value = 1246.92
$ws.Range("J137").Value = 1971
$ws.Range("K137").Value = 3740.76
$ws.Range("L137").Value = 5913
$ws.Range("M137").Value = -1190.76
$ws.Range("N137").Value = -11013

$ws.Range("H138").Value = 538688.4399999999
$ws.Range("I138").Value = 1785
$ws.Range("J138").Value = 605801.4
$ws.Range("K138").Value = 5355
$ws.Range("L138").Value = 1817404.2
$ws.Range("M138").Value = -215
$ws.Range("N138").Value = -1827684.2

$ws.Range("H141").Value = 1754.8667
$ws.Range("I141").Value = 1754.8667
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5264.6001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -84.60009999999966

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3760.6135
$ws.Range("I32").Value = 4227.1665
$ws.Range("J32").Value = 1661.125
$ws.Range("K32").Value = 4227.1665
$ws.Range("L32").Value = 1661.125
$ws.Range("M32").Value = -3940.1665
$ws.Range("N32").Value = -2235.125

$ws.Range("H61").Value = 1949.3
$ws.Range("I61").Value = 1711.2858
$ws.Range("J61").Value = 2504.6667
$ws.Range("K61").Value = 1711.2858
$ws.Range("L61").Value = 2504.6667
$ws.Range("M61").Value = -1499.2858
$ws.Range("N61").Value = -2928.6667

$ws.Range("H74").Value = 1730.5454
$ws.Range("I74").Value = 1579.4286
$ws.Range("J74").Value = 1995
$ws.Range("K74").Value = 1579.4286
$ws.Range("L74").Value = 1995
$ws.Range("M74").Value = -705.4286
$ws.Range("N74").Value = -3743

$ws.Range("H77").Value = 1730.5454
$ws.Range("I77").Value = 1579.4286
$ws.Range("J77").Value = 1995
$ws.Range("K77").Value = 7897.143
$ws.Range("L77").Value = 9975
$ws.Range("M77").Value = -3529.143
$ws.Range("N77").Value = -18711

$ws.Range("H110").Value = 1343.5
$ws.Range("I110").Value = 1168.762
$ws.Range("J110").Value = 5013
$ws.Range("K110").Value = 1168.762
$ws.Range("L110").Value = 5013
$ws.Range("M110").Value = 876.2380000000001
$ws.Range("N110").Value = -9103

$ws.Range("H132").Value = 2527.46
$ws.Range("I132").Value = 2287.5278
$ws.Range("J132").Value = 3144.4285
$ws.Range("K132").Value = 6862.5834
$ws.Range("L132").Value = 9433.2855
$ws.Range("M132").Value = -4332.5834
$ws.Range("N132").Value = -14493.2855

$ws.Range("H136").Value = 1949.3
$ws.Range("I136").Value = 1711.2858
$ws.Range("J136").Value = 2504.6667
$ws.Range("K136").Value = 5133.857400000001
$ws.Range("L136").Value = 7514.000100000001
$ws.Range("M136").Value = -2583.857400000001
$ws.Range("N136").Value = -12614.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2651.068
$ws.Range("I86").Value = 2834.4482
$ws.Range("J86").Value = 2296.5334
$ws.Range("K86").Value = 2834.4482
$ws.Range("L86").Value = 2296.5334
$ws.Range("M86").Value = -1711.4482
$ws.Range("N86").Value = -4542.5334

$ws.Range("H89").Value = 2651.068
$ws.Range("I89").Value = 2834.4482
$ws.Range("J89").Value = 2296.5334
$ws.Range("K89").Value = 14172.241
$ws.Range("L89").Value = 11482.667
$ws.Range("M89").Value = -8556.240999999998
$ws.Range("N89").Value = -22714.667

$ws.Range("H107").Value = 1475.55
$ws.Range("I107").Value = 1034.6428
$ws.Range("J107").Value = 2504.3333
$ws.Range("K107").Value = 1034.6428
$ws.Range("L107").Value = 2504.3333
$ws.Range("M107").Value = 885.3571999999999
$ws.Range("N107").Value = -6344.3333

$ws.Range("H132").Value = 49398.4
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 49398.4
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 49398.4
$ws.Range("N132").Value = -59518.4

$ws.Range("H134").Value = 5497.6523
$ws.Range("I134").Value = 1137.2354
$ws.Range("J134").Value = 17852.166
$ws.Range("K134").Value = 3411.7062
$ws.Range("L134").Value = 53556.49800000001
$ws.Range("M134").Value = -876.7062000000001
$ws.Range("N134").Value = -58626.49800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 7937759
$ws.Range("I134").Value = 10753889
$ws.Range("J134").Value = 1391.1818
$ws.Range("K134").Value = 32261667
$ws.Range("L134").Value = 4173.5454
$ws.Range("M134").Value = -32259132
$ws.Range("N134").Value = -9243.545399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 365.6
$ws.Range("I12").Value = 70.333336
$ws.Range("J12").Value = 439.41666
$ws.Range("K12").Value = 211.000008
$ws.Range("L12").Value = 1318.24998
$ws.Range("M12").Value = -38.00000800000001
$ws.Range("N12").Value = -1664.24998

$ws.Range("H68").Value = 1970.7693
$ws.Range("I68").Value = 600
$ws.Range("J68").Value = 1997.6471
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 5992.9413
$ws.Range("M68").Value = -989
$ws.Range("N68").Value = -7614.9413

$ws.Range("H69").Value = 2454.2666
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 2415.2856
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 7245.8568
$ws.Range("M69").Value = -8189
$ws.Range("N69").Value = -8867.856800000001

$ws.Range("H71").Value = 1970.7693
$ws.Range("I71").Value = 600
$ws.Range("J71").Value = 1997.6471
$ws.Range("K71").Value = 5400
$ws.Range("L71").Value = 17978.8239
$ws.Range("M71").Value = -1344
$ws.Range("N71").Value = -26090.8239

$ws.Range("H72").Value = 2454.2666
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 2415.2856
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 21737.5704
$ws.Range("M72").Value = -22944
$ws.Range("N72").Value = -29849.5704

$ws.Range("H122").Value = 622.6
$ws.Range("I122").Value = 553.5
$ws.Range("J122").Value = 726.25
$ws.Range("K122").Value = 4981.5
$ws.Range("L122").Value = 6536.25
$ws.Range("M122").Value = -2531.5
$ws.Range("N122").Value = -11436.25

$ws.Range("H131").Value = 30304722
$ws.Range("I131").Value = 500000400
$ws.Range("J131").Value = 1774.1936
$ws.Range("K131").Value = 1500001200
$ws.Range("L131").Value = 5322.5808
$ws.Range("M131").Value = -1499996160
$ws.Range("N131").Value = -15402.5808

$ws.Range("H141").Value = 1870.8
$ws.Range("I141").Value = 1870.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5612.4
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -432.3999999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2661.3928
$ws.Range("I132").Value = 2073.7273
$ws.Range("J132").Value = 4816.1665
$ws.Range("K132").Value = 6221.1819
$ws.Range("L132").Value = 14448.4995
$ws.Range("M132").Value = -3691.1819
$ws.Range("N132").Value = -19508.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2800
$ws.Range("I16").Value = 2800
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2800
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2630

$ws.Range("H40").Value = 2504
$ws.Range("I40").Value = 2504
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2504
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2368

$ws.Range("H68").Value = 1484.1666
$ws.Range("I68").Value = 1311.4
$ws.Range("J68").Value = 1700.125
$ws.Range("K68").Value = 1311.4
$ws.Range("L68").Value = 1700.125
$ws.Range("M68").Value = -562.4000000000001
$ws.Range("N68").Value = -3198.125

$ws.Range("H71").Value = 1484.1666
$ws.Range("I71").Value = 1311.4
$ws.Range("J71").Value = 1700.125
$ws.Range("K71").Value = 6557
$ws.Range("L71").Value = 8500.625
$ws.Range("M71").Value = -2813
$ws.Range("N71").Value = -15988.625

$ws.Range("H122").Value = 47229556
$ws.Range("I122").Value = 94449450
$ws.Range("J122").Value = 9666.333000000001
$ws.Range("K122").Value = 283348350
$ws.Range("L122").Value = 28998.999
$ws.Range("M122").Value = -283345900
$ws.Range("N122").Value = -33898.999

$ws.Range("H132").Value = 33261.438
$ws.Range("I132").Value = 1562.0454
$ws.Range("J132").Value = 103000.1
$ws.Range("K132").Value = 4686.1362
$ws.Range("L132").Value = 309000.3
$ws.Range("M132").Value = -2156.1362
$ws.Range("N132").Value = -314060.3

$ws.Range("H136").Value = 6219.8096
$ws.Range("I136").Value = 8841.77
$ws.Range("J136").Value = 1959.125
$ws.Range("K136").Value = 26525.31
$ws.Range("L136").Value = 5877.375
$ws.Range("M136").Value = -23975.31
$ws.Range("N136").Value = -10977.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4000
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3459

$ws.Range("H104").Value = 22500
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 22500
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 22500
$ws.Range("M104").Value = -29488

$ws.Range("H122").Value = 28899224
$ws.Range("I122").Value = 43346170
$ws.Range("J122").Value = 5335
$ws.Range("K122").Value = 130038510
$ws.Range("L122").Value = 16005
$ws.Range("M122").Value = -130036060
$ws.Range("N122").Value = -20905

$ws.Range("H126").Value = 55556852
$ws.Range("I126").Value = 92593360
$ws.Range("J126").Value = 2092.375
$ws.Range("K126").Value = 277780080
$ws.Range("L126").Value = 6277.125
$ws.Range("M126").Value = -277777610
$ws.Range("N126").Value = -11217.125

$ws.Range("H138").Value = 31964.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 31964.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 31964.5
$ws.Range("N138").Value = -42244.5
